{"js": "// Update the cover sheet's academic year from (2023/24) to (2024/25).\nconst searchResults = context.document.body.search(\"(2023/24)\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"(2024/25)\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the cover sheet's academic year from (2023/24) to (2024/25).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"2023/24\"\n$find.Replacement.Text = \"2024/25\"\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, \"wdReplaceAll\")\n"}
